$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "all": append 2020-05-09 (serial 43960) data row, pushing the
# existing footnote row down by one.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Activate()

# Push the trailing footnote row (old row 32) down to row 33 and inherit
# the formatting of the row above (row 31) for the new row 32.
$wsAll.Rows("32:32").Insert()

$wsAll.Range("A32").Value = 43960
$wsAll.Range("B32").Value = 276
$wsAll.Range("C32").Value = 271
$wsAll.Range("D32").Value = 80
$wsAll.Range("E32").Value = 70
$wsAll.Range("F32").Value = 10
$wsAll.Range("G32").Value = 8
$wsAll.Range("H32").Value = 183

$null = $wsAll.Range("A32").Select()

# ---------------------------------------------------------------------
# Sheet "kobe": correct the 2020-05-08 row, then append the 2020-05-09
# data row (pushing the footnote row down by one).
# ---------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Activate()

$wsKobe.Range("D86").Value = 3
$wsKobe.Range("E86").Value = 276

$wsKobe.Rows("87:87").Insert()

$wsKobe.Range("A87").Value = 43960
$wsKobe.Range("B87").Value = 0
$wsKobe.Range("C87").Value = 2480
$wsKobe.Range("D87").Value = 0
$wsKobe.Range("E87").Value = 276
$wsKobe.Range("F87").Value = 75
$wsKobe.Range("G87").Value = 66
$wsKobe.Range("H87").Value = 9
$wsKobe.Range("I87").Value = 8
$wsKobe.Range("J87").Value = 174

$null = $wsKobe.Range("A87").Select()

# ---------------------------------------------------------------------
# Sheet "other": append the 2020-05-09 data row, pushing the footnote
# row down by one.
# ---------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Activate()

$wsOther.Rows("62:62").Insert()

$wsOther.Range("A62").Value = 43960
$wsOther.Range("B62").Value = 0
$wsOther.Range("C62").Value = 14
$wsOther.Range("D62").Value = 5
$wsOther.Range("E62").Value = 4
$wsOther.Range("F62").Value = 1
$wsOther.Range("G62").Value = 0
$wsOther.Range("H62").Value = 9

$null = $wsOther.Range("A62").Select()

# Restore "all" as the active sheet/tab, matching the original workbook.
$wsAll.Activate()
$null = $wsAll.Range("A32").Select()
